$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 10002
$ws.Range("I70").Value = 10002
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 30006
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -29736
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 10002
$ws.Range("I73").Value = 10002
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 30006
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -29070
$ws.Range("N73").ClearContents()
# Row 80
$ws.Range("H80").Value = 572.14703
$ws.Range("I80").Value = 580.0909
$ws.Range("J80").Value = 568.34784
$ws.Range("K80").Value = 1740.2727
$ws.Range("L80").Value = 1705.04352
$ws.Range("M80").Value = -742.2727
$ws.Range("N80").Value = -3701.04352
# Row 83
$ws.Range("H83").Value = 572.14703
$ws.Range("I83").Value = 580.0909
$ws.Range("J83").Value = 568.34784
$ws.Range("K83").Value = 5220.8181
$ws.Range("L83").Value = 5115.130560000001
$ws.Range("M83").Value = -228.8181000000004
$ws.Range("N83").Value = -15099.13056
# Row 86
$ws.Range("H86").Value = 1812.3684
$ws.Range("I86").Value = 1957.9231
$ws.Range("J86").Value = 1497
$ws.Range("K86").Value = 1957.9231
$ws.Range("L86").Value = 1497
$ws.Range("M86").Value = -834.9231
$ws.Range("N86").Value = -3743
# Row 89
$ws.Range("H89").Value = 1812.3684
$ws.Range("I89").Value = 1957.9231
$ws.Range("J89").Value = 1497
$ws.Range("K89").Value = 9789.6155
$ws.Range("L89").Value = 7485
$ws.Range("M89").Value = -4173.6155
$ws.Range("N89").Value = -18717
# Row 121
$ws.Range("H121").Value = 1597.5
$ws.Range("J121").Value = 1833.3334
$ws.Range("L121").Value = 5500.0002
$ws.Range("N121").Value = -8994.0002
# Row 129
$ws.Range("H129").Value = 889.614
$ws.Range("J129").Value = 892.64813
$ws.Range("L129").Value = 2677.94439
$ws.Range("N129").Value = -12677.94439
# Row 131
$ws.Range("H131").Value = 1780
$ws.Range("I131").Value = 1336
$ws.Range("J131").Value = 4000
$ws.Range("K131").Value = 4008
$ws.Range("L131").Value = 12000
$ws.Range("M131").Value = 1032
$ws.Range("N131").Value = -22080
# Row 137
$ws.Range("H137").Value = 1506.6666
$ws.Range("I137").Value = 1436.3636
$ws.Range("J137").Value = 1700
$ws.Range("K137").Value = 4309.0908
$ws.Range("L137").Value = 5100
$ws.Range("M137").Value = -1759.0908
$ws.Range("N137").Value = -10200
# Row 138
$ws.Range("H138").Value = 5447.6626
$ws.Range("I138").Value = 1014.7059
$ws.Range("J138").Value = 6589.485
$ws.Range("K138").Value = 3044.1177
$ws.Range("L138").Value = 19768.455
$ws.Range("M138").Value = 2095.8823
$ws.Range("N138").Value = -30048.455

$ws = $wb.Worksheets.Item("ARM")
# Row 30
$ws.Range("H30").Value = 2411
$ws.Range("I30").Value = 1451.2
$ws.Range("K30").Value = 1451.2
$ws.Range("M30").Value = -1301.2
# Row 45
$ws.Range("H45").Value = 13578.889
$ws.Range("I45").Value = 21282.2
$ws.Range("J45").Value = 3949.75
$ws.Range("K45").Value = 21282.2
$ws.Range("L45").Value = 3949.75
$ws.Range("M45").Value = -20905.2
$ws.Range("N45").Value = -4703.75
# Row 102
$ws.Range("H102").Value = 6175372.5
$ws.Range("I102").Value = 7409447
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 7409447
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -7407825
$ws.Range("N102").Value = -8244

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1669.4333
$ws.Range("I94").Value = 1397.6364
$ws.Range("J94").Value = 2416.875
$ws.Range("K94").Value = 1397.6364
$ws.Range("L94").Value = 2416.875
$ws.Range("M94").Value = -946.6364000000001
$ws.Range("N94").Value = -3318.875

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 1000000000
$ws.Range("J17").Value = 1000000000
$ws.Range("L17").Value = 1000000000
$ws.Range("N17").Value = -1000000348
# Row 22
$ws.Range("H22").Value = 465.13333
$ws.Range("I22").Value = 417.14285
$ws.Range("K22").Value = 417.14285
$ws.Range("M22").Value = -67.14285000000001
# Row 31
$ws.Range("H31").Value = 8034.222
$ws.Range("I31").Value = 1925.1666
$ws.Range("J31").Value = 20252.334
$ws.Range("K31").Value = 1925.1666
$ws.Range("L31").Value = 20252.334
$ws.Range("M31").Value = -1630.1666
$ws.Range("N31").Value = -20842.334
# Row 34
$ws.Range("H34").Value = 8034.222
$ws.Range("I34").Value = 1925.1666
$ws.Range("J34").Value = 20252.334
$ws.Range("K34").Value = 1925.1666
$ws.Range("L34").Value = 20252.334
$ws.Range("M34").Value = -1723.1666
$ws.Range("N34").Value = -20656.334
# Row 58
$ws.Range("H58").Value = 1148.25
$ws.Range("I58").Value = 890.9048
$ws.Range("J58").Value = 1920.2858
$ws.Range("K58").Value = 890.9048
$ws.Range("L58").Value = 1920.2858
$ws.Range("M58").Value = -687.9048
$ws.Range("N58").Value = -2326.2858
# Row 86
$ws.Range("H86").Value = 2369.5
$ws.Range("I86").Value = 2321.6667
$ws.Range("K86").Value = 2321.6667
$ws.Range("M86").Value = -1198.6667
# Row 89
$ws.Range("H89").Value = 2369.5
$ws.Range("I89").Value = 2321.6667
$ws.Range("K89").Value = 11608.3335
$ws.Range("M89").Value = -5992.333500000001
# Row 132
$ws.Range("H132").Value = 3936
$ws.Range("I132").Value = 3299.7144
$ws.Range("J132").Value = 5049.5
$ws.Range("K132").Value = 9899.143199999999
$ws.Range("L132").Value = 15148.5
$ws.Range("M132").Value = -7369.143199999999
$ws.Range("N132").Value = -20208.5
# Row 134
$ws.Range("H134").Value = 5369.9165
$ws.Range("I134").Value = 5993.9
$ws.Range("J134").Value = 2250
$ws.Range("K134").Value = 17981.7
$ws.Range("L134").Value = 6750
$ws.Range("M134").Value = -15446.7
$ws.Range("N134").Value = -11820
# Row 135
$ws.Range("H135").Value = 36200
$ws.Range("J135").Value = 36200
$ws.Range("L135").Value = 36200
$ws.Range("N135").Value = -46340
# Row 136
$ws.Range("H136").Value = 1148.25
$ws.Range("I136").Value = 890.9048
$ws.Range("J136").Value = 1920.2858
$ws.Range("K136").Value = 2672.7144
$ws.Range("L136").Value = 5760.857400000001
$ws.Range("M136").Value = -122.7143999999998
$ws.Range("N136").Value = -10860.8574

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 501072.84
$ws.Range("I5").Value = 574.8
$ws.Range("J5").Value = 858571.4399999999
$ws.Range("K5").Value = 1724.4
$ws.Range("L5").Value = 2575714.32
$ws.Range("M5").Value = -1612.4
$ws.Range("N5").Value = -2575938.32
# Row 131
$ws.Range("H131").Value = 21164342
$ws.Range("J131").Value = 25642058
$ws.Range("L131").Value = 76926174
$ws.Range("N131").Value = -76936254
# Row 135
$ws.Range("H135").Value = 501072.84
$ws.Range("I135").Value = 574.8
$ws.Range("J135").Value = 858571.4399999999
$ws.Range("K135").Value = 5173.2
$ws.Range("L135").Value = 7727142.959999999
$ws.Range("M135").Value = -2638.2
$ws.Range("N135").Value = -7732212.959999999

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1071.375
$ws.Range("I97").Value = 932.3077
$ws.Range("J97").Value = 1674
$ws.Range("K97").Value = 932.3077
$ws.Range("L97").Value = 1674
$ws.Range("M97").Value = -436.3077
$ws.Range("N97").Value = -2666
# Row 102
$ws.Range("H102").Value = 4639.9
$ws.Range("I102").Value = 1199.75
$ws.Range("J102").Value = 6933.3335
$ws.Range("K102").Value = 1199.75
$ws.Range("L102").Value = 6933.3335
$ws.Range("M102").Value = 422.25
$ws.Range("N102").Value = -10177.3335
# Row 122
$ws.Range("H122").Value = 4052882.5
$ws.Range("I122").Value = 5893992
$ws.Range("J122").Value = 2441.2
$ws.Range("K122").Value = 17681976
$ws.Range("L122").Value = 7323.599999999999
$ws.Range("M122").Value = -17679526
$ws.Range("N122").Value = -12223.6
# Row 132
$ws.Range("H132").Value = 3235.8845
$ws.Range("I132").Value = 2278.125
$ws.Range("J132").Value = 3661.5557
$ws.Range("K132").Value = 6834.375
$ws.Range("L132").Value = 10984.6671
$ws.Range("M132").Value = -4304.375
$ws.Range("N132").Value = -16044.6671

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 69506.87
$ws.Range("I7").Value = 93145.73
$ws.Range("K7").Value = 93145.73
$ws.Range("M7").Value = -93033.73
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
# Row 45
$ws.Range("H45").Value = 9000
$ws.Range("I45").Value = 3000
$ws.Range("K45").Value = 3000
$ws.Range("M45").Value = -2593
# Row 55
$ws.Range("H55").Value = 15625314
$ws.Range("I55").Value = 217.91667
$ws.Range("J55").Value = 25000370
$ws.Range("K55").Value = 217.91667
$ws.Range("L55").Value = 25000370
$ws.Range("M55").Value = -44.91667000000001
$ws.Range("N55").Value = -25000716
# Row 122
$ws.Range("H122").Value = 5432940.5
$ws.Range("I122").Value = 5957192
$ws.Range("K122").Value = 17871576
$ws.Range("M122").Value = -17869126
# Row 126
$ws.Range("H126").Value = 69506.87
$ws.Range("I126").Value = 93145.73
$ws.Range("K126").Value = 279437.19
$ws.Range("M126").Value = -276967.19
# Row 132
$ws.Range("H132").Value = 14948512
$ws.Range("I132").Value = 19703708
$ws.Range("J132").Value = 3611.1428
$ws.Range("K132").Value = 59111124
$ws.Range("L132").Value = 10833.4284
$ws.Range("M132").Value = -59108594
$ws.Range("N132").Value = -15893.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 22
$ws.Range("H22").Value = 6117.5
$ws.Range("J22").Value = 6117.5
$ws.Range("L22").Value = 6117.5
$ws.Range("N22").Value = -6703.5
# Row 122
$ws.Range("H122").Value = 3032.6
$ws.Range("I122").Value = 1721.1666
$ws.Range("J122").Value = 4999.75
$ws.Range("K122").Value = 5163.4998
$ws.Range("L122").Value = 14999.25
$ws.Range("M122").Value = -2713.4998
$ws.Range("N122").Value = -19899.25
# Row 126
$ws.Range("H126").Value = 1653.4
$ws.Range("I126").Value = 1237.3334
$ws.Range("J126").Value = 2277.5
$ws.Range("K126").Value = 3712.0002
$ws.Range("L126").Value = 6832.5
$ws.Range("M126").Value = -1242.0002
$ws.Range("N126").Value = -11772.5
# Row 132
$ws.Range("H132").Value = 1888.4
$ws.Range("I132").Value = 1204.2858
$ws.Range("K132").Value = 3612.8574
$ws.Range("M132").Value = -1082.8574
